$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.290.04"
$ws.Range("D3").Value = "'1.550.83"
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'209.84"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").Value = "'0.479"
$ws.Range("E6").Value = "  -2.32%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "'23.79"
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("E10").Value = "  -1.62%  "
$ws.Range("D11").Value = "'0.0890"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "'1.773.10"
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("D13").Value = "'1.545.07"
$ws.Range("E13").Value = "  -2.13%  "
$ws.Range("D14").Value = "'28.296.78"
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("D16").Value = "'3.61"
$ws.Range("E16").Value = "  -1.96%  "
$ws.Range("D17").Value = "'60.78"
$ws.Range("E17").Value = "  -2.56%  "
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("E20").Value = "  -2.63%  "
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("E23").Value = "  -2.83%  "
$ws.Range("D24").Value = "'2.03"
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("D25").Value = "'150.84"
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("E26").Value = "  -1.95%  "
$ws.Range("E27").Value = "  -1.30%  "
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("E29").Value = "  -3.52%  "
$ws.Range("E30").Value = "  -2.75%  "
$ws.Range("E31").Value = "  -4.80%  "
$ws.Range("E32").Value = "  -1.64%  "
$ws.Range("D33").Value = "'1.387.98"
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("E34").Value = "  -3.16%  "
$ws.Range("E35").Value = "  +2.31%  "
$ws.Range("E36").Value = "  -3.92%  "
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("E39").Value = "  -3.04%  "
$ws.Range("D41").Value = "'1.90"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("E43").Value = "  -2.30%  "
$ws.Range("D44").Value = "'0.0458"
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("E45").Value = "  -2.45%  "
$ws.Range("D46").Value = "'61.82"
$ws.Range("E46").Value = "  -2.15%  "
$ws.Range("D47").Value = "'1.685.19"
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("D48").Value = "'0.895"
$ws.Range("E48").Value = "  -6.94%  "
$ws.Range("D49").Value = "'85.44"
$ws.Range("E49").Value = "  -1.37%  "
$ws.Range("E50").Value = "  +7.28%  "
$ws.Range("E51").Value = "  +0.35%  "
